# Updated cryptos list (price and volume columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" '27.285.75'
Set-TextCell "E2" '  +1.21%  '

Set-TextCell "D3" '1.830.97'
Set-TextCell "E3" '  +0.67%  '

Set-TextCell "D4" '1.012'
Set-TextCell "E4" '  +0.94%  '

Set-TextCell "D5" '314.53'
Set-TextCell "E5" '  +1.50%  '

Set-TextCell "E6" '  +0.78%  '

Set-TextCell "D7" '0.4733'
Set-TextCell "E7" '  +1.70%  '

Set-TextCell "D8" '0.3684'

Set-TextCell "D9" '0.07440'
Set-TextCell "E9" '  +1.12%  '

Set-TextCell "D10" '0.8852'
Set-TextCell "E10" '  +1.51%  '

Set-TextCell "E11" '  +0.94%  '

Set-TextCell "D12" '1.879.53'
Set-TextCell "E12" '  +1.46%  '

Set-TextCell "D13" '0.07329'
Set-TextCell "E13" '  +3.04%  '

Set-TextCell "D14" '94.15'
Set-TextCell "E14" '  +2.99%  '

Set-TextCell "D15" '5.427'
Set-TextCell "E15" '  +0.57%  '

Set-TextCell "D16" '6.556'
Set-TextCell "E16" '  +0.70%  '

Set-TextCell "D17" '1.009'
Set-TextCell "E17" '  +0.54%  '

Set-TextCell "E18" '  +1.10%  '

Set-TextCell "D20" '27.582.36'
Set-TextCell "E20" '  +2.19%  '

Set-TextCell "E21" '  +0.76%  '

Set-TextCell "D22" '5.284'
Set-TextCell "E22" '  -0.17%  '

Set-TextCell "E23" '  +0.89%  '

Set-TextCell "D24" '2.107.22'
Set-TextCell "E24" '  +2.40%  '

Set-TextCell "D25" '1.891'
Set-TextCell "E25" '  -0.18%  '

Set-TextCell "D26" '151.88'
Set-TextCell "E26" '  +0.72%  '

Set-TextCell "E27" '  +1.45%  '

Set-TextCell "D28" '2.140'
Set-TextCell "E28" '  -0.05%  '

Set-TextCell "D29" '5.228'
Set-TextCell "E29" '  -0.36%  '

Set-TextCell "E30" '  +0.58%  '

Set-TextCell "D31" '0.08989'
Set-TextCell "E31" '  +1.04%  '

Set-TextCell "D32" '0.7483'
Set-TextCell "E32" '  -1.31%  '

Set-TextCell "D33" '1.174'
Set-TextCell "E33" '  +0.80%  '

Set-TextCell "D34" '4.542'
Set-TextCell "E34" '  +0.88%  '

Set-TextCell "D35" '2.950'
Set-TextCell "E35" '  +1.69%  '

Set-TextCell "E36" '  +0.90%  '

Set-TextCell "D37" '1.093'
Set-TextCell "E37" '  +0.31%  '

Set-TextCell "D38" '0.05339'
Set-TextCell "E38" '  +0.79%  '

Set-TextCell "E39" '  +0.43%  '

Set-TextCell "D40" '2.424'
Set-TextCell "E40" '  +3.49%  '

Set-TextCell "D41" '2.963'
Set-TextCell "E41" '  -0.40%  '

Set-TextCell "D42" '7.239'
Set-TextCell "E42" '  +0.72%  '

Set-TextCell "E43" '  +0.12%  '

Set-TextCell "D44" '0.1657'
Set-TextCell "E44" '  -0.10%  '

Set-TextCell "D45" '8.485'
Set-TextCell "E45" '  +0.49%  '

Set-TextCell "D46" '0.4923'
Set-TextCell "E46" '  +1.19%  '

Set-TextCell "E47" '  +1.40%  '

Set-TextCell "D48" '105.07'
Set-TextCell "E48" '  +1.47%  '

Set-TextCell "E49" '  +0.88%  '

Set-TextCell "D50" '1.666'
Set-TextCell "E50" '  +0.08%  '

Set-TextCell "D51" '0.06301'
Set-TextCell "E51" '  +0.16%  '
